$d = $word.ActiveDocument

# 1. Paragraph 3: drop the trailing period.
$d.Content.Find.Execute(
    "3. Once user successfully bid a ride, his other bids will be deleted.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "3. Once user successfully bid a ride, his other bids will be deleted", 2) | Out-Null

# 2. Drop the (hidden) _GoBack bookmark that lives in the trailing empty paragraph.
try {
    $d.Bookmarks("_GoBack").Delete()
} catch {
}

# 3. Turn that trailing paragraph into constraint 4, and append three more
#    paragraphs: constraints 5 and 6, plus two blank trailing paragraphs
#    (one still styled "p1", the final one back to "Normal").
$last = $d.Paragraphs.Last
$last.Range.InsertAfter("4. Driver can only post rides one hour after current time.`r5. Driver can only have one car.`r6. One ride number should only match only one passenger and driver pair for a single ride`r`r")

$count = $d.Paragraphs.Count
$d.Paragraphs($count - 4).Range.Style = "p1"
$d.Paragraphs($count - 3).Range.Style = "p1"
$d.Paragraphs($count - 2).Range.Style = "p1"
$d.Paragraphs($count - 1).Range.Style = "p1"
$d.Paragraphs($count).Range.Style = "Normal"

# 4. Constraint 6's final "." is typed as its own run, like the rest of the
#    sentence was already entered.
$p6end = $d.Paragraphs($count - 2).Range
$p6end.Collapse(0)
$p6end.MoveEnd(1, -1) | Out-Null
$p6end.InsertAfter(".")
$p6end.Bold = 1
$p6end.Bold = 0
